# Boost MOSFET Dissipation Calcs - Rev II PCB update
#
# Updates the two input-variable tables on Sheet1:
#   - Row 4  (main "Input Variables" block): Vin 5 -> 20, Iin 3 -> 5
#   - Row 38 (PD input-spec playground table): Iin 2 -> 3, L(uH) 68 -> 47,
#     Fsw 150000 -> 200000
# All dependent formulas (rows 8-9, 12-16, 19-21, 24-26, 40-42, 46-62 etc.)
# recalculate automatically - no need to touch them directly.
#
# Also nudges the view (zoom + scroll position + selection) to match where
# the author was working after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Input Variables (row 4): Vin, Iin ---
$ws.Range("A4").Value = 20
$ws.Range("B4").Value = 5

# --- Playing-around-with-PD-specs table (row 38): Iin, L(uH), Fsw ---
$ws.Range("B38").Value = 3
$ws.Range("C38").Value = 47
$ws.Range("F38").Value = 200000

# --- View / selection state ---
$win = $excel.ActiveWindow
$win.Zoom = 115
$win.ScrollRow = 34
$win.ScrollColumn = 1
$ws.Range("C39").Select() | Out-Null
